{"js": "// 1) Style font changes (Normal: Liberation Serif -> Times New Roman,\n//    Heading: Liberation Sans -> Arial). cs/eastAsia fonts are left as-is.\nconst styles = context.document.getStyles();\nconst normalStyle = styles.getByNameOrNullObject(\"Normal\");\nconst headingStyle = styles.getByNameOrNullObject(\"Heading\");\nnormalStyle.load(\"isNullObject\");\nheadingStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!normalStyle.isNullObject) {\n  normalStyle.font.name = \"Times New Roman\";\n}\nif (!headingStyle.isNullObject) {\n  headingStyle.font.name = \"Arial\";\n}\nawait context.sync();\n\n// 2) Shrink the \"Heading 18\" paragraph that is still sized at 19pt (sz/szCs\n//    38 half-points) down to 18pt (sz/szCs 36) -- matches the second\n//    \"Heading 18\" occurrence in the document (the first one is already 18pt).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].font.load(\"size\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text.trim() === \"Heading 18\" && paragraph.font.size === 19) {\n    paragraph.font.size = 18;\n    paragraph.font.sizeBidirectional = 18;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Style font changes (Normal: Liberation Serif -> Times New Roman,\n#    Heading: Liberation Sans -> Arial). cs/eastAsia fonts are left as-is.\n$normalStyle = $d.Styles.Item(\"Normal\")\nif ($normalStyle -ne $null) {\n    $normalStyle.Font.Name = \"Times New Roman\"\n}\n\n$headingStyle = $d.Styles.Item(\"Heading\")\nif ($headingStyle -ne $null) {\n    $headingStyle.Font.Name = \"Arial\"\n}\n\n# 2) Shrink the \"Heading 18\" paragraph that is still sized at 19pt (sz/szCs\n#    38 half-points) down to 18pt (sz/szCs 36) -- matches the second\n#    \"Heading 18\" occurrence in the document (the first one is already 18pt).\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    if ($txt -eq \"Heading 18\" -and $p.Range.Font.Size -eq 19) {\n        $p.Range.Font.Size = 18\n        $p.Range.Font.SizeBi = 18\n    }\n}\n"}
